$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Schedule sheet: add a "Status" column for the existing rows, rename the old
# "Release" stage to "Check-in" (moving its date), and append a new row for
# the "Test Cycle 1 Complete + Report" milestone.
# ---------------------------------------------------------------------------
$schedule = $wb.Worksheets.Item("Schedule")

# Existing "Release" row becomes a second "Check-in" row, dated 10/31/2024.
$schedule.Range("A3").Value = "Check-in"
$schedule.Range("B3").Value = 45596

# Populate the new Status column for the two existing data rows.
$schedule.Range("C2").Value = "Done"
$schedule.Range("C3").Value = "To do"

# New row 4: the test-cycle summary milestone, dated 11/1/2024.
$schedule.Range("A4").Value = "Test Cycle 1 Complete + Report"
$schedule.Range("B4").Value = 45597
$schedule.Range("B4").NumberFormat = $schedule.Range("B3").NumberFormat
$schedule.Range("C4").Value = "In-Progress"

# Move the sheet's selection onto the newly added row.
$schedule.Range("A4").Select() | Out-Null

# ---------------------------------------------------------------------------
# Make RTM the active/visible sheet (was Schedule) and move its selection.
# ---------------------------------------------------------------------------
$rtm = $wb.Worksheets.Item("RTM")
$rtm.Activate() | Out-Null
$rtm.Range("H2").Select() | Out-Null
